{"js": "// Replace each old three-digit-division answer (and the header date) with\n// its updated value. Every source string in this document is unique, so an\n// exact, case-sensitive, whole-text search safely targets a single run.\nconst replacements = [\n  [\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"],\n  [\"980\u00f78=122, 4\", \"763\u00f74=190, 3\"],\n  [\"352\u00f74=88, 0\", \"210\u00f74=52, 2\"],\n  [\"732\u00f78=91, 4\", \"650\u00f77=92, 6\"],\n  [\"240\u00f74=60, 0\", \"794\u00f79=88, 2\"],\n  [\"279\u00f73=93, 0\", \"133\u00f76=22, 1\"],\n  [\"625\u00f78=78, 1\", \"505\u00f77=72, 1\"],\n  [\"407\u00f75=81, 2\", \"524\u00f75=104, 4\"],\n  [\"223\u00f73=74, 1\", \"332\u00f79=36, 8\"],\n  [\"673\u00f77=96, 1\", \"859\u00f76=143, 1\"],\n  [\"905\u00f75=181, 0\", \"446\u00f78=55, 6\"],\n  [\"884\u00f75=176, 4\", \"368\u00f74=92, 0\"],\n  [\"466\u00f78=58, 2\", \"482\u00f73=160, 2\"],\n  [\"791\u00f75=158, 1\", \"661\u00f79=73, 4\"],\n  [\"341\u00f76=56, 5\", \"308\u00f77=44, 0\"],\n  [\"504\u00f73=168, 0\", \"106\u00f72=53, 0\"],\n  [\"559\u00f74=139, 3\", \"906\u00f79=100, 6\"],\n  [\"725\u00f79=80, 5\", \"949\u00f73=316, 1\"],\n  [\"738\u00f77=105, 3\", \"606\u00f74=151, 2\"],\n  [\"213\u00f77=30, 3\", \"458\u00f77=65, 3\"],\n  [\"217\u00f73=72, 1\", \"209\u00f72=104, 1\"],\n  [\"813\u00f74=203, 1\", \"432\u00f77=61, 5\"],\n  [\"477\u00f73=159, 0\", \"566\u00f75=113, 1\"],\n  [\"742\u00f74=185, 2\", \"564\u00f74=141, 0\"],\n  [\"687\u00f76=114, 3\", \"520\u00f73=173, 1\"],\n  [\"347\u00f73=115, 2\", \"181\u00f79=20, 1\"]\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the header date and every three-digit-division answer in the table\n# with its new value. Each source string is unique in the document, so a\n# case-sensitive Find/Replace on the exact old text safely targets one run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-03-01 Saturday\", \"2025-03-02 Sunday\"),\n    @(\"980\u00f78=122, 4\", \"763\u00f74=190, 3\"),\n    @(\"352\u00f74=88, 0\", \"210\u00f74=52, 2\"),\n    @(\"732\u00f78=91, 4\", \"650\u00f77=92, 6\"),\n    @(\"240\u00f74=60, 0\", \"794\u00f79=88, 2\"),\n    @(\"279\u00f73=93, 0\", \"133\u00f76=22, 1\"),\n    @(\"625\u00f78=78, 1\", \"505\u00f77=72, 1\"),\n    @(\"407\u00f75=81, 2\", \"524\u00f75=104, 4\"),\n    @(\"223\u00f73=74, 1\", \"332\u00f79=36, 8\"),\n    @(\"673\u00f77=96, 1\", \"859\u00f76=143, 1\"),\n    @(\"905\u00f75=181, 0\", \"446\u00f78=55, 6\"),\n    @(\"884\u00f75=176, 4\", \"368\u00f74=92, 0\"),\n    @(\"466\u00f78=58, 2\", \"482\u00f73=160, 2\"),\n    @(\"791\u00f75=158, 1\", \"661\u00f79=73, 4\"),\n    @(\"341\u00f76=56, 5\", \"308\u00f77=44, 0\"),\n    @(\"504\u00f73=168, 0\", \"106\u00f72=53, 0\"),\n    @(\"559\u00f74=139, 3\", \"906\u00f79=100, 6\"),\n    @(\"725\u00f79=80, 5\", \"949\u00f73=316, 1\"),\n    @(\"738\u00f77=105, 3\", \"606\u00f74=151, 2\"),\n    @(\"213\u00f77=30, 3\", \"458\u00f77=65, 3\"),\n    @(\"217\u00f73=72, 1\", \"209\u00f72=104, 1\"),\n    @(\"813\u00f74=203, 1\", \"432\u00f77=61, 5\"),\n    @(\"477\u00f73=159, 0\", \"566\u00f75=113, 1\"),\n    @(\"742\u00f74=185, 2\", \"564\u00f74=141, 0\"),\n    @(\"687\u00f76=114, 3\", \"520\u00f73=173, 1\"),\n    @(\"347\u00f73=115, 2\", \"181\u00f79=20, 1\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
